# Update trust_name values in column D per diff (Title Case corrections, mojibake cleanup)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(4, 4).Value = 'Trust Fund: School'
$ws.Cells.Item(21, 4).Value = 'Trust Fund: School'
$ws.Cells.Item(23, 4).Value = 'Sovereign Lands, Common Schools, Strategic Investment and Improvement Fund, ND State University, ND Industrial School, School of Mines, Ellendale, School for the Blind, ND School of Science, University of ND, ND State Treasurer, Mayville, Valley City'
$ws.Cells.Item(27, 4).Value = 'Penitentiary Land Fund, School Lands'
$ws.Cells.Item(28, 4).Value = 'Trust Fund: University, Trust Fund: Swamp, Trust Fund: School'
$ws.Cells.Item(29, 4).Value = 'Perm Cmn Schls (Indmty Selec), Perm Common Schls'
$ws.Cells.Item(30, 4).Value = 'Perm Common Schls, Agriculture & Mechanical Cllge'
$ws.Cells.Item(41, 4).Value = 'Trust Fund: Swamp, Trust Fund: School, Trust Fund: Ind School'
$ws.Cells.Item(43, 4).Value = 'Trust Fund: School, Trust Fund: Ind School'
$ws.Cells.Item(44, 4).Value = 'Perm Common Schls, Perm Cmn Schls (Indmty Selec), County Bond, New Mexico Institute of Mining and Technology, Common Schools, Public Schools'
$ws.Cells.Item(51, 4).Value = 'Trust Fund: School'
$ws.Cells.Item(53, 4).Value = 'Trust Fund: Ind School, Trust Fund: Swamp'
$ws.Cells.Item(57, 4).Value = 'Perm Cmn Schls (Indmty Selec), 02-90-0021 Fed Patent'
$ws.Cells.Item(58, 4).Value = 'Perm Common Schls'
$ws.Cells.Item(65, 4).Value = 'Common Schools, Strategic Investment and Improvement Fund, State Hospital, Valley/Mayville, School for the Deaf, ND State University, ND School of Science, Valley City, Veterans Home, Ellendale, Mayville, ND Industrial School'
$ws.Cells.Item(67, 4).Value = 'Rural Credit, Indemnity, School and Public Lands, Strategic Investment and Improvement Fund, Common Schools, ND School of Science, Capitol Building, Veterans Home, Mayville, School of Mines, University of ND, ND Industrial School, ND State University, Ellendale, School for the Deaf, State Hospital, Valley/Mayville'
$ws.Cells.Item(68, 4).Value = 'County Bond, University, St Chrtbl, Penal & Reform Inst, Perm Common Schls, School of Mines'
$ws.Cells.Item(74, 4).Value = 'Trust Fund: Swamp, Trust Fund: School, Trust Fund: Ind School'
$ws.Cells.Item(80, 4).Value = 'Perm Cmn Schls (Indmty Selec), Saline Lands, University of New Mexico, Common Schools, New Mexico School for the Visually Handicapped'
